$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price text is numeric-looking need a Text format first,
# otherwise Excel would coerce the literal string into a number (dropping
# trailing zeros / exact formatting). Restore the default style afterwards
# so no stray formatting is left behind.
$textFirstCells = @("D5", "D6", "D9", "D14", "D17", "D19", "D21", "D24", "D25", "D29", "D30", "D33", "D34", "D35", "D36", "D39", "D40", "D44", "D48", "D50")
foreach ($cellRef in $textFirstCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '57.075.80'
$ws.Range("E2").Value = '  -1.40%  '
$ws.Range("D3").Value = '2.984.04'
$ws.Range("E3").Value = '  -2.33%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '500.35'
$ws.Range("E5").Value = '  -4.87%  '
$ws.Range("D6").Value = '137.60'
$ws.Range("E6").Value = '  -3.40%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  -4.48%  '
$ws.Range("D9").Value = '7.29'
$ws.Range("E9").Value = '  -5.30%  '
$ws.Range("E10").Value = '  -4.44%  '
$ws.Range("E11").Value = '  -4.28%  '
$ws.Range("D12").Value = '3.494.15'
$ws.Range("E12").Value = '  -2.34%  '
$ws.Range("E13").Value = '  -2.42%  '
$ws.Range("D14").Value = '26.06'
$ws.Range("E14").Value = '  -3.69%  '
$ws.Range("E15").Value = '  -5.98%  '
$ws.Range("D16").Value = '57.113.65'
$ws.Range("E16").Value = '  -1.23%  '
$ws.Range("D17").Value = '6.08'
$ws.Range("E17").Value = '  -3.16%  '
$ws.Range("D18").Value = '2.985.51'
$ws.Range("E18").Value = '  -2.00%  '
$ws.Range("D19").Value = '12.62'
$ws.Range("E19").Value = '  -3.79%  '
$ws.Range("E20").Value = '  -3.45%  '
$ws.Range("D21").Value = '320.43'
$ws.Range("E21").Value = '  -5.22%  '
$ws.Range("E23").Value = '  +0.34%  '
$ws.Range("D24").Value = '0.492'
$ws.Range("E24").Value = '  -2.16%  '
$ws.Range("D25").Value = '62.97'
$ws.Range("E25").Value = '  -3.22%  '
$ws.Range("E26").Value = '  +0.16%  '
$ws.Range("E27").Value = '  -5.26%  '
$ws.Range("D28").Value = '0.0₃0892'
$ws.Range("E28").Value = '  -8.65%  '
$ws.Range("D29").Value = '6.61'
$ws.Range("E29").Value = '  -4.73%  '
$ws.Range("D30").Value = '7.11'
$ws.Range("E30").Value = '  -3.76%  '
$ws.Range("E31").Value = '  -4.29%  '
$ws.Range("E32").Value = '  -6.52%  '
$ws.Range("D33").Value = '20.14'
$ws.Range("E33").Value = '  -4.66%  '
$ws.Range("D34").Value = '154.41'
$ws.Range("E34").Value = '  -1.25%  '
$ws.Range("D35").Value = '4.58'
$ws.Range("E35").Value = '  -3.29%  '
$ws.Range("D36").Value = '5.78'
$ws.Range("E36").Value = '  -3.62%  '
$ws.Range("E37").Value = '  -6.68%  '
$ws.Range("E38").Value = '  -7.00%  '
$ws.Range("D39").Value = '0.0665'
$ws.Range("E39").Value = '  -5.55%  '
$ws.Range("D40").Value = '37.80'
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("D41").Value = '3.013.82'
$ws.Range("E41").Value = '  -2.52%  '
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("E43").Value = '  -3.69%  '
$ws.Range("D44").Value = '0.645'
$ws.Range("E44").Value = '  -2.71%  '
$ws.Range("D45").Value = '2.190.51'
$ws.Range("E45").Value = '  -5.79%  '
$ws.Range("E46").Value = '  -6.41%  '
$ws.Range("E47").Value = '  -1.44%  '
$ws.Range("D48").Value = '0.934'
$ws.Range("E48").Value = '  -9.44%  '
$ws.Range("E49").Value = '  -4.82%  '
$ws.Range("D50").Value = '19.15'
$ws.Range("E50").Value = '  -4.55%  '
$ws.Range("E51").Value = '  -11.35%  '

foreach ($cellRef in $textFirstCells) {
    $ws.Range($cellRef).Style = "Normal"
}
